$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B8 value from 12 to 13
$ws.Range("B8").Value = 13

# Update the active selection to B5 (matches sheetView selection change)
$ws.Range("B5").Select()
